# Updates the cryptocurrency price list (columns D = Price, E = Volume(1h))
# and fixes a row-ordering swap between "Hedera" and "TrustWalletToken"
# (rows 37/38), matching a refreshed data pull from coinranking.com.
#
# NOTE on text formatting: column D/E values are stored as plain text in
# this sheet (not numbers), including values that look numeric (e.g.
# "1.003") or are percentages with padding (e.g. "  -2.36%  "). Assigning
# such a numeric-looking string straight to Range.Value would make Excel
# parse it into a real number, so those cells are written with a leading
# apostrophe (forces text entry, exactly like typing '1.003 into a cell)
# and then their style is reset back to "Normal" so no stray
# quote-prefix/text-format style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value does not get misread as a number/date by Excel
# (two-dot "thousands" prices, percentage strings, coin names, URLs, ...)
# -- these can be written as plain text directly.
$plainTextCells = @(
    @{ Cell = "D2"; Value = "28.660.54" },
    @{ Cell = "E2"; Value = "  -2.36%  " },
    @{ Cell = "D3"; Value = "1.888.08" },
    @{ Cell = "E3"; Value = "  -3.35%  " },
    @{ Cell = "E4"; Value = "  -0.10%  " },
    @{ Cell = "E5"; Value = "  -1.78%  " },
    @{ Cell = "E6"; Value = "  -0.08%  " },
    @{ Cell = "E7"; Value = "  -1.48%  " },
    @{ Cell = "E8"; Value = "  -3.80%  " },
    @{ Cell = "E9"; Value = "  -1.87%  " },
    @{ Cell = "E10"; Value = "  -2.66%  " },
    @{ Cell = "E11"; Value = "  -4.15%  " },
    @{ Cell = "E12"; Value = "  -2.07%  " },
    @{ Cell = "D13"; Value = "1.888.85" },
    @{ Cell = "E13"; Value = "  -3.56%  " },
    @{ Cell = "E14"; Value = "  -3.04%  " },
    @{ Cell = "E15"; Value = "  -3.40%  " },
    @{ Cell = "E16"; Value = "  -1.54%  " },
    @{ Cell = "E17"; Value = "  -0.24%  " },
    @{ Cell = "E18"; Value = "  -6.69%  " },
    @{ Cell = "E19"; Value = "  -4.67%  " },
    @{ Cell = "E20"; Value = "  -3.07%  " },
    @{ Cell = "E21"; Value = "  -0.16%  " },
    @{ Cell = "D22"; Value = "28.648.60" },
    @{ Cell = "E22"; Value = "  -2.68%  " },
    @{ Cell = "E23"; Value = "  -3.68%  " },
    @{ Cell = "E24"; Value = "  -3.50%  " },
    @{ Cell = "D25"; Value = "2.122.15" },
    @{ Cell = "E25"; Value = "  -2.93%  " },
    @{ Cell = "E26"; Value = "  -1.98%  " },
    @{ Cell = "E27"; Value = "  -1.75%  " },
    @{ Cell = "E28"; Value = "  -3.23%  " },
    @{ Cell = "E29"; Value = "  -6.47%  " },
    @{ Cell = "E30"; Value = "  -2.68%  " },
    @{ Cell = "E31"; Value = "  -3.66%  " },
    @{ Cell = "E32"; Value = "  -1.71%  " },
    @{ Cell = "E33"; Value = "  -5.94%  " },
    @{ Cell = "E34"; Value = "  -3.85%  " },
    @{ Cell = "E35"; Value = "  -6.95%  " },
    @{ Cell = "E36"; Value = "  -3.81%  " },
    @{ Cell = "B37"; Value = "TrustWalletToken" },
    @{ Cell = "C37"; Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt" },
    @{ Cell = "E37"; Value = "  -2.75%  " },
    @{ Cell = "B38"; Value = "Hedera" },
    @{ Cell = "C38"; Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar" },
    @{ Cell = "E38"; Value = "  -3.46%  " },
    @{ Cell = "E39"; Value = "  +0.04%  " },
    @{ Cell = "E40"; Value = "  -4.47%  " },
    @{ Cell = "E41"; Value = "  -5.14%  " },
    @{ Cell = "E42"; Value = "  -5.12%  " },
    @{ Cell = "E43"; Value = "  -4.41%  " },
    @{ Cell = "E44"; Value = "  -23.63%  " },
    @{ Cell = "E45"; Value = "  -6.52%  " },
    @{ Cell = "E46"; Value = "  +2.94%  " },
    @{ Cell = "E47"; Value = "  -4.19%  " },
    @{ Cell = "E48"; Value = "  -7.64%  " },
    @{ Cell = "E49"; Value = "  -5.10%  " },
    @{ Cell = "E50"; Value = "  -2.91%  " },
    @{ Cell = "E51"; Value = "  -2.98%  " }
)

foreach ($item in $plainTextCells) {
    $ws.Range($item.Cell).Value = $item.Value
}

# Cells whose new value looks like a plain number (e.g. "1.003", "322.09")
# and must be forced to stay text via the leading-apostrophe trick, then
# have the resulting style reset to "Normal" to avoid leaving a
# quote-prefix style behind.
$forcedTextCells = @(
    @{ Cell = "D4"; Value = "1.003" },
    @{ Cell = "D5"; Value = "322.09" },
    @{ Cell = "D6"; Value = "1.002" },
    @{ Cell = "D8"; Value = "0.3785" },
    @{ Cell = "D9"; Value = "45.41" },
    @{ Cell = "D10"; Value = "0.07672" },
    @{ Cell = "D11"; Value = "0.9574" },
    @{ Cell = "D12"; Value = "21.86" },
    @{ Cell = "D14"; Value = "6.937" },
    @{ Cell = "D15"; Value = "5.629" },
    @{ Cell = "D16"; Value = "0.07015" },
    @{ Cell = "D17"; Value = "1.004" },
    @{ Cell = "D18"; Value = "82.42" },
    @{ Cell = "D19"; Value = "0.000009455" },
    @{ Cell = "D20"; Value = "16.57" },
    @{ Cell = "D23"; Value = "5.319" },
    @{ Cell = "D24"; Value = "10.82" },
    @{ Cell = "D26"; Value = "2.078" },
    @{ Cell = "D27"; Value = "155.09" },
    @{ Cell = "D28"; Value = "18.92" },
    @{ Cell = "D29"; Value = "5.587" },
    @{ Cell = "D31"; Value = "1.810" },
    @{ Cell = "D32"; Value = "0.09214" },
    @{ Cell = "D33"; Value = "0.8394" },
    @{ Cell = "D34"; Value = "5.035" },
    @{ Cell = "D35"; Value = "1.239" },
    @{ Cell = "D36"; Value = "3.053" },
    @{ Cell = "D37"; Value = "1.144" },
    @{ Cell = "D38"; Value = "0.05605" },
    @{ Cell = "D39"; Value = "1.002" },
    @{ Cell = "D40"; Value = "0.02021" },
    @{ Cell = "D41"; Value = "7.417" },
    @{ Cell = "D42"; Value = "0.5452" },
    @{ Cell = "D43"; Value = "0.1741" },
    @{ Cell = "D44"; Value = "0.000002938" },
    @{ Cell = "D45"; Value = "9.138" },
    @{ Cell = "D46"; Value = "2.692" },
    @{ Cell = "D47"; Value = "0.5139" },
    @{ Cell = "D48"; Value = "11.15" },
    @{ Cell = "D49"; Value = "2.080" },
    @{ Cell = "D50"; Value = "0.06755" },
    @{ Cell = "D51"; Value = "110.45" }
)

foreach ($item in $forcedTextCells) {
    $ws.Range($item.Cell).Value = "'" + $item.Value
    $ws.Range($item.Cell).Style = "Normal"
}
